# narcotics_bs.xlsx - "get order list service implemented and put together all services"
#
# The sheet's "packaging" columns are collapsed down to a single running
# "quantity" column:
#   - the old UnitsPerPackage column (G) is dropped
#   - the old PackagesInStorage column (H) becomes the new "quantity" column (G)
#   - the old ThresholdPack column (I) is dropped
#   - the two trailing blank columns (J, K) are dropped
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "UnitsPerPackage" column (G). This shifts PackagesInStorage
# (old H) into G and ThresholdPack (old I) into H.
$ws.Columns.Item(7).Delete()

# Drop the now-shifted "ThresholdPack" column (H).
$ws.Columns.Item(8).Delete()

# The two trailing empty columns (originally J and K) have already been
# pulled in by the two deletes above, so the sheet now ends at column I
# (H and I are the blank header-styled cells with no values).

# Rename the header for the surviving quantity column.
$ws.Range("G1").Value = "quantity"

# The old PackagesInStorage numbers were right-aligned (style index 3);
# the merged "quantity" column goes back to the plain/default look.
$ws.Range("G2:G7").ClearFormats()

# Restate the explicit row heights (matches the resave in the source file).
for ($r = 1; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# Selection / scroll position left behind by the editing session.
[void]$ws.Range("G16").Select()
